$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Column C ("Förändrad") holds a date serial number that gets bumped by
# one day for every data row (rows 2 through 490).
$ws.Range("C2:C490").Value = 46061
